# Update "想去人数" (column F) figures across the three data sheets
# that report per-event interest counts ("展览", "演出", "全部类型").
# The numbers below are new scrape totals replacing the prior ones.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        [string]$SheetName,
        [hashtable]$RowToValue
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowToValue.Keys) {
        $ws.Cells.Item($row, 6).Value = $RowToValue[$row]
    }
}

# NOTE: this runtime's PowerShell parser does not reliably bind named
# (-Param value) arguments, so Set-FValues is always invoked positionally
# below: Set-FValues <SheetName> <RowToValue>.

# Sheet "展览" (sheet1.xml)
$exhibitionUpdates = @{
    2  = 1296
    3  = 1178
    4  = 888
    5  = 107
    7  = 652
    8  = 101
    9  = 49
    10 = 38
    11 = 2347
    12 = 1589
    13 = 1344
    15 = 232
    16 = 546
    17 = 756
    18 = 42
    19 = 284
    20 = 1087
    22 = 15
    24 = 4624
    26 = 115
    27 = 41
    29 = 126
    30 = 202
    31 = 82
    32 = 14
    33 = 672
    34 = 452
    35 = 64
    36 = 39
    37 = 236
    38 = 371
    39 = 968
    40 = 126
    41 = 92
    42 = 145
    43 = 116
}
Set-FValues "展览" $exhibitionUpdates

# Sheet "演出" (sheet2.xml)
$performanceUpdates = @{
    3  = 781
    5  = 427
    11 = 16
}
Set-FValues "演出" $performanceUpdates

# Sheet "全部类型" (sheet4.xml)
$allTypesUpdates = @{
    2  = 1296
    4  = 781
    5  = 1178
    6  = 888
    8  = 427
    9  = 107
    11 = 652
    12 = 101
    13 = 49
    17 = 38
    18 = 2347
    19 = 1589
    20 = 1344
    22 = 232
    23 = 546
    25 = 756
    26 = 42
    27 = 284
    28 = 1087
    29 = 15
    30 = 4624
    32 = 41
    34 = 126
    35 = 202
    36 = 82
    37 = 14
    38 = 672
    39 = 452
    40 = 64
    41 = 371
    42 = 968
    43 = 126
    44 = 92
    45 = 145
    46 = 116
    48 = 16
}
Set-FValues "全部类型" $allTypesUpdates
